# Euclid_Macro_Path_Plot.xlsx - "Excel Plotter, Updated Qidi Prts"
# Update the variable_position_* X/Y coordinates on Sheet1.
# All downstream values (R/S helper columns, chart caches, series
# names, data-label text) are driven off these source cells by
# formulas, so updating them here and letting Excel recalculate
# reproduces the rest of the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# variable_position_preflight (row 6)
$ws.Range("C6").Value = 35
$ws.Range("D6").Value = 300

# variable_position_side (row 7)
$ws.Range("C7").Value = 35
$ws.Range("D7").Value = 320

# variable_position_dock (row 8)
$ws.Range("C8").Value = 35
$ws.Range("D8").Value = 349

# variable_position_zstop (row 9)
$ws.Range("C9").Value = 235
$ws.Range("D9").Value = 349

# variable_position_exit (row 12)
$ws.Range("C12").Value = 70
$ws.Range("D12").Value = 349

$excel.CalculateFullRebuild()
$wb.RefreshAll()

# Match the saved selection/active-cell state captured in the diff.
$ws.Range("C12").Select()
